$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force the cell to Text format so numeric-looking strings (e.g. "1.00",
    # "14.30", percentages) are not auto-converted by Excel into numbers,
    # then restore the default "Normal" style so no stray number format lingers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" '63.797.71'
Set-TextCell $ws "E2" '  -0.94%  '

# Row 3
Set-TextCell $ws "D3" '3.079.70'
Set-TextCell $ws "E3" '  -7.27%  '

# Row 4
Set-TextCell $ws "E4" '  -1.41%  '

# Row 5
Set-TextCell $ws "D5" '587.57'
Set-TextCell $ws "E5" '  -1.18%  '

# Row 6
Set-TextCell $ws "D6" '155.68'
Set-TextCell $ws "E6" '  +4.88%  '

# Row 7
Set-TextCell $ws "D7" '1.00'
Set-TextCell $ws "E7" '  -1.10%  '

# Row 8
Set-TextCell $ws "D8" '0.538'
Set-TextCell $ws "E8" '  +0.50%  '

# Row 9
Set-TextCell $ws "D9" '3.078.44'
Set-TextCell $ws "E9" '  -3.09%  '

# Row 10
Set-TextCell $ws "E10" '  -4.30%  '

# Row 11
Set-TextCell $ws "E11" '  -3.47%  '

# Row 12
Set-TextCell $ws "E12" '  -3.22%  '

# Row 13
Set-TextCell $ws "E13" '  -4.54%  '

# Row 14
Set-TextCell $ws "D14" '36.82'
Set-TextCell $ws "E14" '  -2.98%  '

# Row 15
Set-TextCell $ws "E15" '  -2.04%  '

# Row 16
Set-TextCell $ws "D16" '3.583.71'
Set-TextCell $ws "E16" '  -7.27%  '

# Row 17
Set-TextCell $ws "D17" '7.17'
Set-TextCell $ws "E17" '  -2.67%  '

# Row 18
Set-TextCell $ws "D18" '63.687.32'
Set-TextCell $ws "E18" '  -1.03%  '

# Row 19
Set-TextCell $ws "D19" '3.075.37'
Set-TextCell $ws "E19" '  -4.70%  '

# Row 20
Set-TextCell $ws "D20" '471.35'
Set-TextCell $ws "E20" '  -0.73%  '

# Row 21
Set-TextCell $ws "D21" '14.30'
Set-TextCell $ws "E21" '  -2.44%  '

# Row 22
Set-TextCell $ws "D22" '0.705'
Set-TextCell $ws "E22" '  -5.23%  '

# Row 23
Set-TextCell $ws "E23" '  -2.93%  '

# Row 24
Set-TextCell $ws "D24" '2.44'
Set-TextCell $ws "E24" '  -2.37%  '

# Row 25
Set-TextCell $ws "E25" '  -5.74%  '

# Row 26
Set-TextCell $ws "D26" '80.43'
Set-TextCell $ws "E26" '  -2.98%  '

# Row 27
Set-TextCell $ws "D27" '10.37'
Set-TextCell $ws "E27" '  +3.06%  '

# Row 28
Set-TextCell $ws "E28" '  -0.21%  '

# Row 29
Set-TextCell $ws "D29" '7.37'
Set-TextCell $ws "E29" '  +1.11%  '

# Row 30
Set-TextCell $ws "E30" '  -3.27%  '

# Row 31
Set-TextCell $ws "E31" '  -0.94%  '

# Row 32
Set-TextCell $ws "E32" '  -6.71%  '

# Row 33
Set-TextCell $ws "E33" '  -7.40%  '

# Row 34
Set-TextCell $ws "D34" '27.13'
Set-TextCell $ws "E34" '  -4.65%  '

# Row 35
Set-TextCell $ws "D35" '0.0₃0829'
Set-TextCell $ws "E35" '  -3.61%  '

# Row 36
Set-TextCell $ws "E36" '  -3.09%  '

# Row 37
Set-TextCell $ws "E37" '  -4.61%  '

# Row 38
Set-TextCell $ws "D38" '3.26'
Set-TextCell $ws "E38" '  -2.75%  '

# Row 39
Set-TextCell $ws "E39" '  -5.48%  '

# Row 40
Set-TextCell $ws "D40" '50.61'
Set-TextCell $ws "E40" '  -2.19%  '

# Row 41
Set-TextCell $ws "D41" '9.12'
Set-TextCell $ws "E41" '  -3.66%  '

# Row 42
Set-TextCell $ws "D42" '432.81'
Set-TextCell $ws "E42" '  -7.90%  '

# Row 43
Set-TextCell $ws "D43" '0.290'
Set-TextCell $ws "E43" '  -3.22%  '

# Row 44
Set-TextCell $ws "D44" '0.111'
Set-TextCell $ws "E44" '  +0.82%  '

# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws "D45" '0.0360'
Set-TextCell $ws "E45" '  -4.70%  '

# Row 46
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell $ws "D46" '39.90'
Set-TextCell $ws "E46" '  +1.06%  '

# Row 47
Set-TextCell $ws "D47" '2.809.52'
Set-TextCell $ws "E47" '  -4.58%  '

# Row 48
Set-TextCell $ws "D48" '129.94'
Set-TextCell $ws "E48" '  -2.22%  '

# Row 49
Set-TextCell $ws "E49" '  +0.06%  '

# Row 50
Set-TextCell $ws "D50" '24.95'
Set-TextCell $ws "E50" '  -0.07%  '

# Row 51
Set-TextCell $ws "E51" '  -3.90%  '
